$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7188875675201416
$ws.Range("B1").Value = 3.642741441726685
$ws.Range("C1").Value = 2.669186353683472
$ws.Range("D1").Value = 2.207858085632324
$ws.Range("E1").Value = 2.004126787185669
